$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price/Volume(1h) columns) per latest scrape.
# D-column price values are prefixed with a leading apostrophe so Excel
# keeps numeric-looking strings (e.g. '1.00', '7.70') as text instead of
# coercing them to Number cells and dropping the trailing zeros.

$ws.Range('D2').Value = '''76.102.38'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').Value = '''2.855.75'
$ws.Range('E3').Value = '  +7.64%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''193.98'
$ws.Range('E5').Value = '  +4.44%  '
$ws.Range('D6').Value = '''599.13'
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '''0.551'
$ws.Range('E8').Value = '  +3.69%  '
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('D10').Value = '''2.850.50'
$ws.Range('E10').Value = '  +7.48%  '
$ws.Range('D11').Value = '''0.392'
$ws.Range('E11').Value = '  +10.70%  '
$ws.Range('E12').Value = '  -2.00%  '
$ws.Range('D13').Value = '''4.91'
$ws.Range('E13').Value = '  +4.13%  '
$ws.Range('D14').Value = '''3.380.77'
$ws.Range('E14').Value = '  +7.73%  '
$ws.Range('D15').Value = '''75.932.14'
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').Value = '''27.51'
$ws.Range('E16').Value = '  +4.34%  '
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('D18').Value = '''2.843.89'
$ws.Range('E18').Value = '  +7.37%  '
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('D20').Value = '''12.52'
$ws.Range('E20').Value = '  +5.73%  '
$ws.Range('D21').Value = '''382.83'
$ws.Range('E21').Value = '  +3.15%  '
$ws.Range('D22').Value = '''2.33'
$ws.Range('E22').Value = '  +2.60%  '
$ws.Range('D23').Value = '''4.14'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('D24').Value = '''72.04'
$ws.Range('E24').Value = '  +3.72%  '
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  +2.95%  '
$ws.Range('D27').Value = '''2.988.33'
$ws.Range('E27').Value = '  +7.08%  '
$ws.Range('E28').Value = '  +4.64%  '
$ws.Range('D29').Value = '''0.0000105'
$ws.Range('E29').Value = '  +12.22%  '
$ws.Range('D30').Value = '''0.997'
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('D31').Value = '''1.42'
$ws.Range('E31').Value = '  +1.35%  '
$ws.Range('D32').Value = '''518.97'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').Value = '''7.70'
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('E34').Value = '  +4.80%  '
$ws.Range('D35').Value = '''0.998'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').Value = '''166.43'
$ws.Range('E36').Value = '  +2.64%  '
$ws.Range('E37').Value = '  +4.82%  '
$ws.Range('E38').Value = '  +0.89%  '
$ws.Range('D39').Value = '''19.43'
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('D40').Value = '''186.66'
$ws.Range('E40').Value = '  +10.84%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = '''5.10'
$ws.Range('E42').Value = '  +2.92%  '
$ws.Range('E43').Value = '  +5.79%  '
$ws.Range('D44').Value = '''1.68'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('E45').Value = '  +5.01%  '
$ws.Range('D46').Value = '''40.30'
$ws.Range('E46').Value = '  +3.23%  '
$ws.Range('D47').Value = '''0.0888'
$ws.Range('E47').Value = '  +5.10%  '
$ws.Range('D48').Value = '''2.38'
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('E49').Value = '  +9.79%  '
$ws.Range('D50').Value = '''3.76'
$ws.Range('E50').Value = '  +4.00%  '
$ws.Range('D51').Value = '''0.662'
$ws.Range('E51').Value = '  +12.39%  '
